# Hapus referensi: 118976 - TEST
# Row 79 (A79:F79 = 015 | Kementerian Keuangan | 118976 | TEST | TESTT TESTTTT | 78)
# is the dummy "TEST" reference row that must be removed. Deleting it shifts
# the used range up to A1:F78 and the "No" column (F) in the remaining rows
# needs to stay a contiguous 1..77 numbering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row for the "118976 - TEST" reference entry.
$ws.Rows.Item(79).Delete()

# Renumber the "No" column (F) sequentially as plain numbers (1..77) for the
# remaining data rows (rows 2..78), matching the post-delete sequence.
for ($r = 2; $r -le 78; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 1
}
